$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ57338166"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9906.54288309394
$ws.Cells.Item(2, 3).Value = [double]"1.800380683246211e-15"
$ws.Cells.Item(3, 2).Value = -1031.213237402815
$ws.Cells.Item(3, 3).Value = 0.1565465266498265
$ws.Cells.Item(4, 2).Value = -501.7252305353898
$ws.Cells.Item(4, 3).Value = 0.3337839314360908
$ws.Cells.Item(5, 2).Value = 264.6610658113797
$ws.Cells.Item(5, 3).Value = 0.2112866245778373
$ws.Cells.Item(6, 2).Value = 132.1766929505734
$ws.Cells.Item(6, 3).Value = 0.579535029560079
$ws.Cells.Item(7, 2).Value = 42.17772245926076
$ws.Cells.Item(7, 3).Value = 0.8735298866443663
$ws.Cells.Item(8, 2).Value = -39.71145313775801
$ws.Cells.Item(8, 3).Value = 0.8704822724988371
$ws.Cells.Item(9, 2).Value = 17.28323433600566
$ws.Cells.Item(9, 3).Value = 0.8249385545871765
$ws.Cells.Item(10, 2).Value = -1289.309234329043
$ws.Cells.Item(10, 3).Value = [double]"5.569666296341531e-13"
$ws.Cells.Item(11, 2).Value = -35.42839637265212
$ws.Cells.Item(11, 3).Value = [double]"1.438772577683192e-05"
$ws.Cells.Item(12, 2).Value = 175.293511247277
$ws.Cells.Item(12, 3).Value = 0.01230700331675299
$ws.Cells.Item(13, 2).Value = 453.5547601998289
$ws.Cells.Item(13, 3).Value = [double]"1.20621213945493e-22"
$ws.Cells.Item(14, 2).Value = 0.2236423132203778
$ws.Cells.Item(14, 3).Value = [double]"4.512662846812229e-06"
$ws.Cells.Item(15, 2).Value = [double]"2.685831699777234e-06"
$ws.Cells.Item(15, 3).Value = 0.9568383093336469
$ws.Cells.Item(16, 2).Value = -22.29551922430231
$ws.Cells.Item(16, 3).Value = 0.04881426408264395
$ws.Cells.Item(17, 2).Value = -3.338786103322133
$ws.Cells.Item(17, 3).Value = 0.6083524194282345
$ws.Cells.Item(18, 2).Value = -2733.349492082581
$ws.Cells.Item(18, 3).Value = 0.005218353126404175
$ws.Cells.Item(19, 2).Value = -3870.506354889485
$ws.Cells.Item(19, 3).Value = 0.0002222431746986263

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ57748283"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 10107.29081152513
$ws.Cells.Item(2, 3).Value = [double]"2.372339833086613e-16"
$ws.Cells.Item(3, 2).Value = -712.1665564532349
$ws.Cells.Item(3, 3).Value = 0.3114156671322352
$ws.Cells.Item(4, 2).Value = -527.0858010042947
$ws.Cells.Item(4, 3).Value = 0.306894824382467
$ws.Cells.Item(5, 2).Value = 99.06804833998612
$ws.Cells.Item(5, 3).Value = 0.6338279493372083
$ws.Cells.Item(6, 2).Value = 195.7473573445334
$ws.Cells.Item(6, 3).Value = 0.4059372187579549
$ws.Cells.Item(7, 2).Value = 56.97973830877191
$ws.Cells.Item(7, 3).Value = 0.8266188677865605
$ws.Cells.Item(8, 2).Value = 108.6008564380531
$ws.Cells.Item(8, 3).Value = 0.6503887073685743
$ws.Cells.Item(9, 2).Value = -85.89257668243448
$ws.Cells.Item(9, 3).Value = 0.2676945863610501
$ws.Cells.Item(10, 2).Value = -1403.996198921379
$ws.Cells.Item(10, 3).Value = [double]"1.50428968734867e-15"
$ws.Cells.Item(11, 2).Value = -34.30180275614915
$ws.Cells.Item(11, 3).Value = [double]"2.288185701949101e-05"
$ws.Cells.Item(12, 2).Value = 209.2329819781789
$ws.Cells.Item(12, 3).Value = 0.002262963617954393
$ws.Cells.Item(13, 2).Value = 471.7762118493035
$ws.Cells.Item(13, 3).Value = [double]"1.345982948868512e-24"
$ws.Cells.Item(14, 2).Value = 0.2235922430074882
$ws.Cells.Item(14, 3).Value = [double]"3.658797169832912e-06"
$ws.Cells.Item(15, 2).Value = [double]"-1.41171573251875e-05"
$ws.Cells.Item(15, 3).Value = 0.7746172185733877
$ws.Cells.Item(16, 2).Value = -20.17131938778724
$ws.Cells.Item(16, 3).Value = 0.0714324407001386
$ws.Cells.Item(17, 2).Value = -4.537892893797303
$ws.Cells.Item(17, 3).Value = 0.4818276970094425
$ws.Cells.Item(18, 2).Value = -2305.033464868777
$ws.Cells.Item(18, 3).Value = 0.01689269135476067
$ws.Cells.Item(19, 2).Value = -4055.091938674924
$ws.Cells.Item(19, 3).Value = [double]"9.07147497965258e-05"

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ58229681"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 8969.108888469353
$ws.Cells.Item(2, 3).Value = [double]"5.333480207436166e-13"
$ws.Cells.Item(3, 2).Value = -936.0130970957825
$ws.Cells.Item(3, 3).Value = 0.1848134385485178
$ws.Cells.Item(4, 2).Value = -328.1292203041307
$ws.Cells.Item(4, 3).Value = 0.5307814236384878
$ws.Cells.Item(5, 2).Value = 65.06255022110723
$ws.Cells.Item(5, 3).Value = 0.7570529602725682
$ws.Cells.Item(6, 2).Value = -42.23782722439971
$ws.Cells.Item(6, 3).Value = 0.8584348462544077
$ws.Cells.Item(7, 2).Value = 14.51106917757193
$ws.Cells.Item(7, 3).Value = 0.9559654652604792
$ws.Cells.Item(8, 2).Value = -125.9122059816884
$ws.Cells.Item(8, 3).Value = 0.6022878063696724
$ws.Cells.Item(9, 2).Value = 8.913944475300028
$ws.Cells.Item(9, 3).Value = 0.908797128237835
$ws.Cells.Item(10, 2).Value = -1450.790995144369
$ws.Cells.Item(10, 3).Value = [double]"3.123061396114603e-16"
$ws.Cells.Item(11, 2).Value = -30.64686065294207
$ws.Cells.Item(11, 3).Value = 0.0001603849546635446
$ws.Cells.Item(12, 2).Value = 217.8915443313504
$ws.Cells.Item(12, 3).Value = 0.00165642247610735
$ws.Cells.Item(13, 2).Value = 484.5780901071357
$ws.Cells.Item(13, 3).Value = [double]"8.776565555031562e-26"
$ws.Cells.Item(14, 2).Value = 0.2184824323214948
$ws.Cells.Item(14, 3).Value = [double]"6.244770390076333e-06"
$ws.Cells.Item(15, 2).Value = [double]"-5.113422642993325e-05"
$ws.Cells.Item(15, 3).Value = 0.2928146349297963
$ws.Cells.Item(16, 2).Value = -13.44825820724765
$ws.Cells.Item(16, 3).Value = 0.2321509800214407
$ws.Cells.Item(17, 2).Value = 1.23789724732151
$ws.Cells.Item(17, 3).Value = 0.8486991821745938
$ws.Cells.Item(18, 2).Value = -1957.492859553397
$ws.Cells.Item(18, 3).Value = 0.04253510846892501
$ws.Cells.Item(19, 2).Value = -3858.802746896961
$ws.Cells.Item(19, 3).Value = 0.0002239205324387199

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ58748508"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 8295.331349633241
$ws.Cells.Item(2, 3).Value = [double]"1.185870296012957e-11"
$ws.Cells.Item(3, 2).Value = -1000.902549551033
$ws.Cells.Item(3, 3).Value = 0.1483498510712902
$ws.Cells.Item(4, 2).Value = 66.32321789103094
$ws.Cells.Item(4, 3).Value = 0.8972142399988372
$ws.Cells.Item(5, 2).Value = 217.7198357380515
$ws.Cells.Item(5, 3).Value = 0.2963226954467588
$ws.Cells.Item(6, 2).Value = 136.3623896554375
$ws.Cells.Item(6, 3).Value = 0.5604289887734122
$ws.Cells.Item(7, 2).Value = 169.4288855717078
$ws.Cells.Item(7, 3).Value = 0.5166364195107134
$ws.Cells.Item(8, 2).Value = 106.2488540290401
$ws.Cells.Item(8, 3).Value = 0.6577579538529728
$ws.Cells.Item(9, 2).Value = -54.06695357926819
$ws.Cells.Item(9, 3).Value = 0.4849077456251112
$ws.Cells.Item(10, 2).Value = -1189.844073321672
$ws.Cells.Item(10, 3).Value = [double]"1.249022716437592e-11"
$ws.Cells.Item(11, 2).Value = -30.84477353700022
$ws.Cells.Item(11, 3).Value = 0.0001178861247411925
$ws.Cells.Item(12, 2).Value = 224.6557980875934
$ws.Cells.Item(12, 3).Value = 0.001089557330146548
$ws.Cells.Item(13, 2).Value = 501.6510839593374
$ws.Cells.Item(13, 3).Value = [double]"4.245366958842465e-28"
$ws.Cells.Item(14, 2).Value = 0.2217475585995917
$ws.Cells.Item(14, 3).Value = [double]"3.949749433790856e-06"
$ws.Cells.Item(15, 2).Value = [double]"-1.725454227195865e-05"
$ws.Cells.Item(15, 3).Value = 0.7245042616297341
$ws.Cells.Item(16, 2).Value = -13.17196736519356
$ws.Cells.Item(16, 3).Value = 0.2374253723419607
$ws.Cells.Item(17, 2).Value = 3.540620006049329
$ws.Cells.Item(17, 3).Value = 0.5798191612858355
$ws.Cells.Item(18, 2).Value = -2081.886894807644
$ws.Cells.Item(18, 3).Value = 0.03053376990749051
$ws.Cells.Item(19, 2).Value = -4074.154201505211
$ws.Cells.Item(19, 3).Value = [double]"7.67894187723053e-05"

# --- Worksheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ59311285"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9898.012420113602
$ws.Cells.Item(2, 3).Value = [double]"1.903424877402149e-15"
$ws.Cells.Item(3, 2).Value = -1060.52237947245
$ws.Cells.Item(3, 3).Value = 0.1446045965405017
$ws.Cells.Item(4, 2).Value = -377.0448234299421
$ws.Cells.Item(4, 3).Value = 0.4667016373821661
$ws.Cells.Item(5, 2).Value = 75.4777110156575
$ws.Cells.Item(5, 3).Value = 0.7201809382755745
$ws.Cells.Item(6, 2).Value = 34.29075736180852
$ws.Cells.Item(6, 3).Value = 0.8856832694343945
$ws.Cells.Item(7, 2).Value = 147.3231700727268
$ws.Cells.Item(7, 3).Value = 0.5767073214835066
$ws.Cells.Item(8, 2).Value = 186.4066013685762
$ws.Cells.Item(8, 3).Value = 0.4461697242886079
$ws.Cells.Item(9, 2).Value = -4.070091729929366
$ws.Cells.Item(9, 3).Value = 0.9585638653803787
$ws.Cells.Item(10, 2).Value = -1261.814498163843
$ws.Cells.Item(10, 3).Value = [double]"1.592546458771345e-12"
$ws.Cells.Item(11, 2).Value = -32.37091134987151
$ws.Cells.Item(11, 3).Value = [double]"7.553025208049137e-05"
$ws.Cells.Item(12, 2).Value = 200.6230055732732
$ws.Cells.Item(12, 3).Value = 0.003969194072922422
$ws.Cells.Item(13, 2).Value = 454.075974623673
$ws.Cells.Item(13, 3).Value = [double]"8.389183080632302e-23"
$ws.Cells.Item(14, 2).Value = 0.2386943292209807
$ws.Cells.Item(14, 3).Value = [double]"8.205356418659653e-07"
$ws.Cells.Item(15, 2).Value = [double]"-2.182928514226495e-05"
$ws.Cells.Item(15, 3).Value = 0.6559311096093668
$ws.Cells.Item(16, 2).Value = -19.21732583882019
$ws.Cells.Item(16, 3).Value = 0.09034757414780867
$ws.Cells.Item(17, 2).Value = -3.359479059629828
$ws.Cells.Item(17, 3).Value = 0.6068760358824858
$ws.Cells.Item(18, 2).Value = -2633.916397633797
$ws.Cells.Item(18, 3).Value = 0.006952665716928307
$ws.Cells.Item(19, 2).Value = -4439.745350471936
$ws.Cells.Item(19, 3).Value = [double]"2.245207349627664e-05"

# --- Worksheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ59858130"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9345.868301811199
$ws.Cells.Item(2, 3).Value = [double]"5.267986212770827e-14"
$ws.Cells.Item(3, 2).Value = -911.4153297962168
$ws.Cells.Item(3, 3).Value = 0.1969755345743665
$ws.Cells.Item(4, 2).Value = -466.3293530061959
$ws.Cells.Item(4, 3).Value = 0.387011802683043
$ws.Cells.Item(5, 2).Value = 56.20963388997699
$ws.Cells.Item(5, 3).Value = 0.7893410105244477
$ws.Cells.Item(6, 2).Value = 67.87061920179437
$ws.Cells.Item(6, 3).Value = 0.7749436565158203
$ws.Cells.Item(7, 2).Value = 211.8124699428391
$ws.Cells.Item(7, 3).Value = 0.4226653964914737
$ws.Cells.Item(8, 2).Value = 110.517643897941
$ws.Cells.Item(8, 3).Value = 0.6484664730110621
$ws.Cells.Item(9, 2).Value = -78.30099890842973
$ws.Cells.Item(9, 3).Value = 0.3157652774213826
$ws.Cells.Item(10, 2).Value = -1410.890240343264
$ws.Cells.Item(10, 3).Value = [double]"2.296832362123221e-15"
$ws.Cells.Item(11, 2).Value = -33.59837455797924
$ws.Cells.Item(11, 3).Value = [double]"4.032738908149722e-05"
$ws.Cells.Item(12, 2).Value = 231.5961534030031
$ws.Cells.Item(12, 3).Value = 0.0008365328052759458
$ws.Cells.Item(13, 2).Value = 448.7172761025232
$ws.Cells.Item(13, 3).Value = [double]"3.596862905988984e-22"
$ws.Cells.Item(14, 2).Value = 0.2408254995449589
$ws.Cells.Item(14, 3).Value = [double]"7.214059183774334e-07"
$ws.Cells.Item(15, 2).Value = [double]"-2.515977820008418e-05"
$ws.Cells.Item(15, 3).Value = 0.6071451001671972
$ws.Cells.Item(16, 2).Value = -16.42594218465917
$ws.Cells.Item(16, 3).Value = 0.1482818426562431
$ws.Cells.Item(17, 2).Value = 2.18879248098463
$ws.Cells.Item(17, 3).Value = 0.7348131504359359
$ws.Cells.Item(18, 2).Value = -2559.598375996118
$ws.Cells.Item(18, 3).Value = 0.008766588357667573
$ws.Cells.Item(19, 2).Value = -4284.141468648579
$ws.Cells.Item(19, 3).Value = [double]"4.782633539633276e-05"

# --- Worksheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ00425709"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9636.370790633675
$ws.Cells.Item(2, 3).Value = [double]"8.201876903229636e-15"
$ws.Cells.Item(3, 2).Value = -827.8635740514894
$ws.Cells.Item(3, 3).Value = 0.2453666330613655
$ws.Cells.Item(4, 2).Value = -167.9738494620709
$ws.Cells.Item(4, 3).Value = 0.7520447868920102
$ws.Cells.Item(5, 2).Value = 97.1725461523379
$ws.Cells.Item(5, 3).Value = 0.6439087946104778
$ws.Cells.Item(6, 2).Value = 51.3173818248405
$ws.Cells.Item(6, 3).Value = 0.8284120318926
$ws.Cells.Item(7, 2).Value = 79.19115059573852
$ws.Cells.Item(7, 3).Value = 0.7638925326195021
$ws.Cells.Item(8, 2).Value = -30.1218002345642
$ws.Cells.Item(8, 3).Value = 0.9010209927628736
$ws.Cells.Item(9, 2).Value = -48.20228094057747
$ws.Cells.Item(9, 3).Value = 0.5396809112728427
$ws.Cells.Item(10, 2).Value = -1456.258624077665
$ws.Cells.Item(10, 3).Value = [double]"2.467721567123044e-16"
$ws.Cells.Item(11, 2).Value = -28.82833117802435
$ws.Cells.Item(11, 3).Value = 0.0004070556586288834
$ws.Cells.Item(12, 2).Value = 182.4165456510785
$ws.Cells.Item(12, 3).Value = 0.00879682310105168
$ws.Cells.Item(13, 2).Value = 472.135668747817
$ws.Cells.Item(13, 3).Value = [double]"2.115594833473749e-24"
$ws.Cells.Item(14, 2).Value = 0.2410784207598202
$ws.Cells.Item(14, 3).Value = [double]"6.587350233784648e-07"
$ws.Cells.Item(15, 2).Value = [double]"-2.814839449803016e-05"
$ws.Cells.Item(15, 3).Value = 0.569911045449131
$ws.Cells.Item(16, 2).Value = -25.13467628364738
$ws.Cells.Item(16, 3).Value = 0.02468590499361624
$ws.Cells.Item(17, 2).Value = -1.698092542765206
$ws.Cells.Item(17, 3).Value = 0.7940273838312915
$ws.Cells.Item(18, 2).Value = -1710.298714523918
$ws.Cells.Item(18, 3).Value = 0.08064590476127852
$ws.Cells.Item(19, 2).Value = -3792.531041336465
$ws.Cells.Item(19, 3).Value = 0.0002633533737886832

# --- Worksheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ00933101"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9940.838854489317
$ws.Cells.Item(2, 3).Value = [double]"6.129973966453716e-16"
$ws.Cells.Item(3, 2).Value = -935.5349939494974
$ws.Cells.Item(3, 3).Value = 0.203722456603545
$ws.Cells.Item(4, 2).Value = -359.9831815174227
$ws.Cells.Item(4, 3).Value = 0.4889917229450595
$ws.Cells.Item(5, 2).Value = 105.2999131244343
$ws.Cells.Item(5, 3).Value = 0.6148972307177459
$ws.Cells.Item(6, 2).Value = 193.6776340343925
$ws.Cells.Item(6, 3).Value = 0.4124456251883216
$ws.Cells.Item(7, 2).Value = 162.4741998530672
$ws.Cells.Item(7, 3).Value = 0.5365796117509114
$ws.Cells.Item(8, 2).Value = 63.53167601230513
$ws.Cells.Item(8, 3).Value = 0.7925080233616387
$ws.Cells.Item(9, 2).Value = -10.16158130272978
$ws.Cells.Item(9, 3).Value = 0.89633778438276
$ws.Cells.Item(10, 2).Value = -1243.293544093795
$ws.Cells.Item(10, 3).Value = [double]"2.299098331311863e-12"
$ws.Cells.Item(11, 2).Value = -33.3233688556751
$ws.Cells.Item(11, 3).Value = [double]"3.979451330235886e-05"
$ws.Cells.Item(12, 2).Value = 205.4276309123536
$ws.Cells.Item(12, 3).Value = 0.003146834436208059
$ws.Cells.Item(13, 2).Value = 461.7472858368235
$ws.Cells.Item(13, 3).Value = [double]"1.121601360270272e-23"
$ws.Cells.Item(14, 2).Value = 0.2227442107574488
$ws.Cells.Item(14, 3).Value = [double]"3.79115903113425e-06"
$ws.Cells.Item(15, 2).Value = [double]"-4.50096850487186e-05"
$ws.Cells.Item(15, 3).Value = 0.3646730760780409
$ws.Cells.Item(16, 2).Value = -23.99760525005607
$ws.Cells.Item(16, 3).Value = 0.03209778643173049
$ws.Cells.Item(17, 2).Value = -4.730233135301233
$ws.Cells.Item(17, 3).Value = 0.4612486001929058
$ws.Cells.Item(18, 2).Value = -1791.994914555151
$ws.Cells.Item(18, 3).Value = 0.06814684565627499
$ws.Cells.Item(19, 2).Value = -3549.930735429057
$ws.Cells.Item(19, 3).Value = 0.0006450847321225601

# --- Worksheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ01446852"
$ws.Rows.Item(20).Delete()
$ws.Cells.Item(2, 2).Value = 9814.36403866569
$ws.Cells.Item(2, 3).Value = [double]"9.725814449297078e-16"
$ws.Cells.Item(3, 2).Value = -616.1244671888005
$ws.Cells.Item(3, 3).Value = 0.3894676835089768
$ws.Cells.Item(4, 2).Value = -137.277856568915
$ws.Cells.Item(4, 3).Value = 0.7872557609561435
$ws.Cells.Item(5, 2).Value = 65.44103895731918
$ws.Cells.Item(5, 3).Value = 0.7537740963398508
$ws.Cells.Item(6, 2).Value = 192.6475131949718
$ws.Cells.Item(6, 3).Value = 0.4126087410531664
$ws.Cells.Item(7, 2).Value = 163.9602752116463
$ws.Cells.Item(7, 3).Value = 0.5297717738519547
$ws.Cells.Item(8, 2).Value = 19.5560511026957
$ws.Cells.Item(8, 3).Value = 0.9350545071905005
$ws.Cells.Item(9, 2).Value = -17.4684231225902
$ws.Cells.Item(9, 3).Value = 0.8213766360081661
$ws.Cells.Item(10, 2).Value = -1335.45403046967
$ws.Cells.Item(10, 3).Value = [double]"3.576645383019023e-14"
$ws.Cells.Item(11, 2).Value = -30.99365077321741
$ws.Cells.Item(11, 3).Value = 0.0001060118105488199
$ws.Cells.Item(12, 2).Value = 188.1086102035151
$ws.Cells.Item(12, 3).Value = 0.006116717437045414
$ws.Cells.Item(13, 2).Value = 491.5024088349977
$ws.Cells.Item(13, 3).Value = [double]"6.04311688116832e-27"
$ws.Cells.Item(14, 2).Value = 0.2249410378863136
$ws.Cells.Item(14, 3).Value = [double]"2.776712544430456e-06"
$ws.Cells.Item(15, 2).Value = [double]"-3.323462996749118e-05"
$ws.Cells.Item(15, 3).Value = 0.4957430211651461
$ws.Cells.Item(16, 2).Value = -23.26598579780688
$ws.Cells.Item(16, 3).Value = 0.03634623392441508
$ws.Cells.Item(17, 2).Value = -4.405201070986301
$ws.Cells.Item(17, 3).Value = 0.4915191293034307
$ws.Cells.Item(18, 2).Value = -2214.701321796439
$ws.Cells.Item(18, 3).Value = 0.0224634346729572
$ws.Cells.Item(19, 2).Value = -3823.840697995059
$ws.Cells.Item(19, 3).Value = 0.0002239226583113833
